$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "2021 - Høst"
$ws.Range("B11").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/hjemme-21-h.pdf)"
$ws.Range("D11").Value = "[Materiale](tidligere-eksamensoppgaver/hjemme-21-h-ekstra.zip)"
$ws.Range("C11").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/hjemme-21-h-solprop.html)"

$ws.Range("C12").Select()
